$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting from column J (rows 3 and 5-12) into the new columns K and L,
# tiling the single-column source across the two-column destination so each row
# of K/L picks up the same cell style as the corresponding row in J.
$ws.Range("J3:J12").Copy()
$ws.Range("K3:L12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The pasted formats for the data rows (5-12) keep J's explicit right alignment;
# the real edit's new styles drop that explicit horizontal alignment (falls back
# to default/general). Clear it back to General for those rows only.
$ws.Range("K5:L12").HorizontalAlignment = 1

# Header years
$ws.Range("K4").Value = 2021
$ws.Range("L4").Value = 2022

# Data values
$ws.Range("K5").Value = 272.60000000000002
$ws.Range("L5").Value = 292.19961890663211

$ws.Range("K7").Value = 98.1
$ws.Range("L7").Value = 99.522498012012946

$ws.Range("K8").Value = 174.5
$ws.Range("L8").Value = 192.67712089461918

$ws.Range("K10").Value = 75.599999999999994
$ws.Range("L10").Value = 88.011952928467494

$ws.Range("K11").Value = 55.5
$ws.Range("L11").Value = 56.919430260413804

$ws.Range("K12").Value = 24.9
$ws.Range("L12").Value = 24.176373211436804

# Match the author's final selection recorded in the saved file
$ws.Range("N5").Select()
